$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.486.95"
$ws.Range("E2").Value = "  +2.36%  "
$ws.Range("D3").Value = "1.827.07"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5073"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3927"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07734"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.267"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("B14").Value = "BinanceUSD"
$ws.Range("C14").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.002"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.563"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "1.824.73"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.00%  "
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06619"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.096"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("D23").Value = "28.494.33"
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.256"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.437"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "2.034.99"
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.130"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1094"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.652"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.659"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07140"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2228"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.007"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02320"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.143"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6240"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5891"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.714"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.972"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.179"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06929"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.02%  "
